$d = $word.ActiveDocument

# Locate the target paragraph (the "Nucleic Acids Research" abstract paragraph)
# robustly via its distinctive leading text, rather than a hard-coded index.
$findRange = $d.Content
$found = $findRange.Find.Execute("While RAG binds and cuts specific regions of the DNA", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$targetParagraph = $findRange.Paragraphs(1).Range

$xml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00A35286" w:rsidRDefault="00A35286" w:rsidP="00A35286"><w:pPr><w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="CMR10" w:eastAsia="Times New Roman" w:hAnsi="CMR10" w:cs="Times New Roman"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">While RAG binds and cuts specific regions of the DNA because of the recognizable sequence patterns, these sites can still vary in sequence within the genome. In a study that we recently published in </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="CMTI10" w:eastAsia="Times New Roman" w:hAnsi="CMTI10" w:cs="Times New Roman"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>Nucleic Acids Research</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="CMR10" w:eastAsia="Times New Roman" w:hAnsi="CMR10" w:cs="Times New Roman"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">, we examined the extent to which RAG will bind and cut the DNA if we modify a binding site sequence at multiple positions and compared the contributions of each nucleotide modification to the collective effect. We illustrate some of our findings in this visual, which is modified from a page in </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="CMR10" w:eastAsia="Times New Roman" w:hAnsi="CMR10" w:cs="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>the Supplementary website that accompanies our publication</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="CMR10" w:eastAsia="Times New Roman" w:hAnsi="CMR10" w:cs="Times New Roman"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">. </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="CMR10" w:hAnsi="CMR10"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>In this interactive visual, w</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="CMR10" w:hAnsi="CMR10"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>e</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="CMR10" w:hAnsi="CMR10"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> show three example comparisons between effects of several single nucleotide changes</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="CMR10" w:hAnsi="CMR10"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> from some starting sequence</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="CMR10" w:hAnsi="CMR10"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> and that of combining these replacements in a single sequence. Through the dropdown menu, one can select any of these three binding site sequences to reveal the effects of the sequence and the individual effects of its constitutive changes. The upper left</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="CMR10" w:hAnsi="CMR10"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> plot</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="CMR10" w:hAnsi="CMR10"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> shows the frequency that RAG creates a DNA loop for the combination of changes to the far right and the individual changes to the left, with position along the sequence where the change was made as the x-axis. The plot on the upper right shows full posterior distributions of the probability that RAG cuts the DNA with the altered sequence. In the bottom row, we present three cumulative distribution functions to show (from left to right) how much time it takes before DNA unloops without cutting, is cut, or a compilation of the two possible fates. </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="CMR10" w:hAnsi="CMR10"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>To</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="CMR10" w:hAnsi="CMR10"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> more easily compare one particular single nucleotide change against the combined changes, hovering the mouse over a colored nucleotide in the sequence below the dropdown menu will send the rest of the data into the background and present only the individual change </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="CMR10" w:hAnsi="CMR10"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>with</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="CMR10" w:hAnsi="CMR10"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> the </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="CMR10" w:hAnsi="CMR10"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>superposition</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="CMR10" w:hAnsi="CMR10"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> of changes. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$targetParagraph.InsertXML($xml) | Out-Null
